# ENH: Add dynamic lapse
# Inserts a new "DynLapse" worksheet between "Lapse" and "ADB" and
# populates it with the dynamic-lapse formula parameter table.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet right after "Lapse" -----------------------------
$lapseSheet = $wb.Worksheets.Item("Lapse")
$dyn = $wb.Worksheets.Add($null, $lapseSheet)
$dyn.Name = "DynLapse"

# --- Header row -------------------------------------------------------------
$dyn.Range("B1").Value = "formula_id"
$dyn.Range("C1").Value = "U"
$dyn.Range("D1").Value = "L"
$dyn.Range("E1").Value = "M"
$dyn.Range("F1").Value = "D"
$dyn.Range("G1").Value = "FactorCap"
$dyn.Range("H1").Value = "FactorFloor"
$dyn.Range("I1").Value = "Y"
$dyn.Range("J1").Value = "Power"

# --- DL001A / DL001B (DL001 formula) ----------------------------------------
$dyn.Range("A2").Value = "DL001A"
$dyn.Range("B2").Value = "DL001"
$dyn.Range("C2").Value = 2
$dyn.Range("D2").Value = 0.5
$dyn.Range("E2").Value = 3.5
$dyn.Range("F2").Value = 0.8
$dyn.Range("C2:D2").NumberFormat = "0%"

$dyn.Range("A3").Value = "DL001B"
$dyn.Range("B3").Value = "DL001"
$dyn.Range("C3").Value = 2
$dyn.Range("D3").Value = 0.5
$dyn.Range("E3").Value = 0.7
$dyn.Range("F3").Value = 1.1
$dyn.Range("C3:D3").NumberFormat = "0%"

# --- DL002A / DL002B (DL002 formula) ----------------------------------------
$dyn.Range("A4").Value = "DL002A"
$dyn.Range("B4").Value = "DL002"
$dyn.Range("G4").Value = 1.5
$dyn.Range("H4").Value = 0.8
$dyn.Range("I4").Value = 1
$dyn.Range("J4").Value = 1

$dyn.Range("A5").Value = "DL002B"
$dyn.Range("B5").Value = "DL002"
$dyn.Range("G5").Value = 2
$dyn.Range("H5").Value = 0.5
$dyn.Range("I5").Value = 1
$dyn.Range("J5").Value = 1

# --- Cosmetics: column widths ------------------------------------------------
$dyn.Columns("A:B").ColumnWidth = 12
$dyn.Range("H1").ColumnWidth = 9.3

# --- Cosmetics: selection / active cell on the new sheet --------------------
$selectResult = $dyn.Range("B8").Select()

# --- Cosmetics: the Inflation sheet view scrolls down a bit -----------------
$inflation = $wb.Worksheets.Item("Inflation")
$inflation.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Restore the originally-active sheet selection (the new DynLapse tab becomes
# active, matching the saved workbook view).
$dyn.Activate()
